$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $startP = $p.Start
    $p.Select()
    $sel = $word.Selection
    $sel.Text = $newText

    # After the assignment above only the paragraph's first run is
    # updated in place (preserving its formatting/xml:space), while any
    # further runs that used to make up the paragraph are still present
    # immediately after it. Remove that now-stale leftover text, but keep
    # the trailing paragraph mark.
    $p2 = $d.Paragraphs.Item($paraIndex).Range
    $delStart = $startP + $newText.Length
    $delEnd = $p2.End - 1
    if ($delEnd -gt $delStart) {
        $d.Range($delStart, $delEnd).Delete()
    }
}

# Title: "Questions: Trigonometric identities (degrees)"
Set-ParagraphText 1 "Questions: Trigonometric identities (degrees)"

# Author: "Dzhemma Ruseva"
Set-ParagraphText 2 "Dzhemma Ruseva"

# Abstract: "A selection of questions on trigonometric identities, where angles are measured in degrees."
Set-ParagraphText 4 "A selection of questions on trigonometric identities, where angles are measured in degrees."
